$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New "Beobachtung" (observation) header column, mirroring the existing
# J1/K1 header-value pair formatting.
$ws.Range("J1").Copy()
$ws.Range("M1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("M1").Value = "Beobachtung"

$ws.Range("K1").Copy()
$ws.Range("N1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("N1").Value = 20

# Duplicate the ESPPRC result value into the new column for row 2
$ws.Range("N2").Value = $ws.Range("H2").Text

# Observation / comment text at the bottom of the table
$ws.Range("M7").Value = "Nur eine verwendete Variable ist nicht elementar. DSSR sinnvoll für ESPPRC Lösung?"

$excel.CutCopyMode = 0

# Update the selection to reflect where the user was working
$ws.Range("N3").Select()
